# Add working set of sequences
#
# The sheet already has columns A:F populated for every data row, but a
# subset of rows ("short" rows) are missing the remaining columns G:N
# (which, for every other row, are populated). Those rows need G:N filled
# in with the same "N/A" value that column F already holds for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5,7,10,11,16,20,27,28,30,33,34,38,39,40,43,47,49,51,54,55,59,62,64,65,70,76,78,80,84,86,88,89,95,101,108,109,110,112,113,117,121,122,126,131,133,139,144,147,149,155,161,162,163,169,170,177,179,180,184,185,186,188,190,193)

foreach ($r in $rows) {
    $ws.Range("G$r`:N$r").Value = "N/A"
}
